$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newValues = @(
    @(8.306412923626578, 2.768804307875526),
    @(8.143944171844646, 2.67820982832475),
    @(7.984307446591386, 2.590934204523033),
    @(7.827425236737846, 2.506822461438918),
    @(7.6732234826847, 2.42572871388097),
    @(7.521631410869367, 2.347515535876427),
    @(7.372581394358958, 2.272053385682949),
    @(7.226008845902345, 2.199220083535496),
    @(7.0818521540604, 2.128900340791164),
    @(6.940052679600399, 2.060985341214664),
    @(6.80055483953812, 1.995372378068071),
    @(6.663306322079896, 1.931964554922573),
    @(6.528258500424876, 1.870670564449233),
    @(6.395367151828792, 1.811404569015091),
    @(6.264593647989086, 1.754086221436944),
    @(6.135906874935531, 1.698640886281588),
    @(6.009286281592359, 1.645000155296232),
    @(5.884726668557478, 1.593102799775229),
    @(5.762245639853067, 1.54289637351257),
    @(5.641895074048588, 1.494339776369626),
    @(5.523778518061655, 1.447407205267492),
    @(5.408076947829414, 1.40209402351133),
    @(5.295085468694369, 1.358425067884943),
    @(5.185262218194373, 1.316465537261784),
    @(5.079285922845465, 1.276333385740656),
    @(4.978107328566682, 1.238209437054657),
    @(4.882961136195439, 1.202337164188827),
    @(4.795288275142211, 1.1690006242884),
    @(4.716530548507789, 1.138472891019122),
    @(4.647828287166298, 1.110944322298286),
    @(4.589746077739643, 1.086461631928708),
    @(4.54216963218478, 1.064910583622269),
    @(4.504411169355498, 1.046048091461703),
    @(4.475434161816517, 1.029559971497696),
    @(4.454077629883037, 1.015115366810553),
    @(4.439213108688921, 1.002402960026531),
    @(4.416246074593184, 0.9881098523062376),
    @(4.390456412385222, 0.9734496117958185),
    @(4.364813616184787, 0.9590846062468814),
    @(4.339317145500344, 0.9450068450200747),
    @(4.313966447116654, 0.931208616337956),
    @(4.288760956023772, 0.9176824753063966),
    @(4.26370009629011, 0.904421232546387),
    @(4.238783281884887, 0.8914179434006844),
    @(4.214009917451979, 0.8786658976814765),
    @(4.189379399038577, 0.8661586099278071),
    @(4.164891114781254, 0.8538898101434369),
    @(4.140544445552172, 0.841853434987786),
    @(4.116338765567718, 0.8300436193943137),
    @(4.092273442961891, 0.8184546885923781),
    @(4.092273442961891, 0.8184546885923781),
    @(3.794242382517293, 0.6865771930269385),
    @(3.518024410562919, 0.5812388156582214),
    @(3.262185848632251, 0.495852248992102),
    @(3.0253289707144, 0.4257870403227675),
    @(2.806119997760789, 0.3676984824652068),
    @(2.603300110974081, 0.3191142071516614),
    @(2.41568812816966, 0.2781701480922638),
    @(2.242179161499242, 0.2434365946770605),
    @(2.08174144710694, 0.2138004729461182),
    @(1.933412458234173, 0.1883837779817913),
    @(1.796294854622382, 0.1664858645747574),
    @(1.669552528995717, 0.1475418513996215),
    @(1.552406859332579, 0.1310921347880844),
    @(1.444133198331093, 0.116759705396982),
    @(1.344057594413594, 0.1042330379341154),
    @(1.251553723123703, 0.09325302250725627),
    @(1.166040003556778, 0.08360286817954254),
    @(1.086976875836521, 0.07510022051234148),
    @(1.013864219432036, 0.0675909479621357)
)

for ($i = 0; $i -lt $newValues.Count; $i++) {
    $r = 2 + $i
    $ws.Cells.Item($r, 3).Value = $newValues[$i][0]
    $ws.Cells.Item($r, 4).Value = $newValues[$i][1]
}

Write-Host "Updated $($newValues.Count) rows (C2:D71)"
